$wb = $excel.ActiveWorkbook

# --- Rename Sheet2 -> Case1_GoogleSearch and populate it (Sample Page Object) ---
$ws2 = $wb.Worksheets.Item("Sheet2")
$ws2.Name = "Case1_GoogleSearch"

# Header row values (write in this order so the new shared-string table
# indices line up: 20=Google Search Keywords, 21=No of Result listed,
# 22=Times(Frequency), 23=Date and Time)
$ws2.Range("A1").Value = "ID"
$ws2.Range("B1").Value = "Google Search Keywords"
$ws2.Range("D1").Value = "No of Result listed"
$ws2.Range("C1").Value = "Times(Frequency)"
$ws2.Range("E1").Value = "Date and Time"

# Header style ("Heading 3" built-in cell style -> bold font + bottom border)
$ws2.Range("A1:E1").Style = "Heading 3"
$ws2.Rows.Item(1).RowHeight = 15

# Column widths
$ws2.Columns.Item(1).ColumnWidth = 12
$ws2.Columns.Item(2).ColumnWidth = 26
$ws2.Columns.Item(3).ColumnWidth = 16.83333333
$ws2.Columns.Item(4).ColumnWidth = 17
$ws2.Columns.Item(5).ColumnWidth = 17

# Freeze the header row, then make this sheet active with F1 selected
$ws2.Activate()
$ws2.Range("A2").Select()
$ws2.Application.ActiveWindow.FreezePanes = $true
$ws2.Range("F1").Select()
